$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of worksheet row number -> corrected K (strikeouts) value.
# These values replace the previous "Strike#"-derived figures in column G ("K")
# with the correct K values, as part of regenerating this save_data sheet.
$gValues = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 1
    6 = 2
    7 = 3
    8 = 3
    9 = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 3
    16 = 1
    17 = 1
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 3
    32 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 2
    39 = 3
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 3
    45 = 1
    46 = 0
    47 = 1
    48 = 1
    49 = 2
    50 = 0
    51 = 2
    52 = 1
    53 = 2
    54 = 0
    55 = 5
    56 = 2
    57 = 0
    58 = 0
    59 = 0
    60 = 2
    61 = 2
    62 = 1
    63 = 0
    64 = 1
    65 = 2
    68 = 1
    69 = 1
    70 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
